$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 600, shifting existing rows 600:701 down to 601:702
$ws.Rows("600:600").Insert()

# Fill in the values for the newly inserted row 600
$ws.Range("A600").Value2 = 6
$ws.Range("B600").Value2 = 'Mercado Mayorista Lo Valledor de Santiago'
$ws.Range("C600").Value2 = 'Metropolitana'
$ws.Range("D600").Value2 = 44694
$ws.Range("E600").Value2 = 13
$ws.Range("F600").Value2 = 100112003
$ws.Range("G600").Value2 = 'Ajo'
$ws.Range("H600").Value2 = 'Chino'
$ws.Range("I600").Value2 = 'Primera'
$ws.Range("J600").Value2 = 3700
$ws.Range("K600").Value2 = 16000
$ws.Range("L600").Value2 = 17000
$ws.Range("M600").Value2 = 16405
$ws.Range("N600").Value2 = '$/caja 10 kilos'
$ws.Range("O600").Value2 = 'China'
$ws.Range("P600").Value2 = 1640
$ws.Range("Q600").Value2 = 10
$ws.Range("R600").Value2 = 'Hortaliza'
